$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top (shifts all existing data down by one row)
$ws.Rows.Item(1).Insert()

# New header row
$ws.Range("A1").Value = "Raum"
$ws.Range("B1").Value = "Kapazität"

# Fix the room-number / capacity bug:
#  - room "008" capacity changes 20 -> 15
#  - room 110 capacity changes 25 -> 20
$ws.Range("B2").Value = 15
$ws.Range("B10").Value = 20

# Update selection to reflect the new last cell
$ws.Range("B15").Select()
